# Updates the cryptocurrency price/volume table on Sheet1 to reflect the
# latest scrape (prices, hourly volume %, and two reordered coin rows).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell reference, new text value.
# Price cells (column D) are written with a leading apostrophe so that
# values such as '0.0000251' or '136.50' are kept as exact text instead
# of being re-interpreted as numbers, then the style is reset to "Normal"
# so no formatting is left behind on the cell.
$updates = @(
    @{ Cell = 'D2'; Value = '67.289.69'; IsPrice = $true }
    @{ Cell = 'E2'; Value = '  +1.45%  '; IsPrice = $false }
    @{ Cell = 'D3'; Value = '3.141.88'; IsPrice = $true }
    @{ Cell = 'E3'; Value = '  +3.58%  '; IsPrice = $false }
    @{ Cell = 'E4'; Value = '  +0.01%  '; IsPrice = $false }
    @{ Cell = 'D5'; Value = '581.26'; IsPrice = $true }
    @{ Cell = 'E5'; Value = '  +0.71%  '; IsPrice = $false }
    @{ Cell = 'D6'; Value = '174.77'; IsPrice = $true }
    @{ Cell = 'E6'; Value = '  +3.90%  '; IsPrice = $false }
    @{ Cell = 'E7'; Value = '  +0.05%  '; IsPrice = $false }
    @{ Cell = 'D8'; Value = '3.136.15'; IsPrice = $true }
    @{ Cell = 'E8'; Value = '  +3.43%  '; IsPrice = $false }
    @{ Cell = 'D9'; Value = '0.526'; IsPrice = $true }
    @{ Cell = 'E9'; Value = '  +1.07%  '; IsPrice = $false }
    @{ Cell = 'D10'; Value = '6.49'; IsPrice = $true }
    @{ Cell = 'E10'; Value = '  -2.95%  '; IsPrice = $false }
    @{ Cell = 'D11'; Value = '0.156'; IsPrice = $true }
    @{ Cell = 'E11'; Value = '  +2.02%  '; IsPrice = $false }
    @{ Cell = 'D12'; Value = '0.485'; IsPrice = $true }
    @{ Cell = 'E12'; Value = '  +0.38%  '; IsPrice = $false }
    @{ Cell = 'D13'; Value = '0.0000251'; IsPrice = $true }
    @{ Cell = 'E13'; Value = '  +1.30%  '; IsPrice = $false }
    @{ Cell = 'D14'; Value = '37.32'; IsPrice = $true }
    @{ Cell = 'E14'; Value = '  +2.09%  '; IsPrice = $false }
    @{ Cell = 'E15'; Value = '  +0.22%  '; IsPrice = $false }
    @{ Cell = 'D16'; Value = '3.647.74'; IsPrice = $true }
    @{ Cell = 'E16'; Value = '  +3.43%  '; IsPrice = $false }
    @{ Cell = 'D17'; Value = '67.223.12'; IsPrice = $true }
    @{ Cell = 'E17'; Value = '  +1.40%  '; IsPrice = $false }
    @{ Cell = 'D18'; Value = '7.20'; IsPrice = $true }
    @{ Cell = 'E18'; Value = '  -0.65%  '; IsPrice = $false }
    @{ Cell = 'D19'; Value = '3.132.28'; IsPrice = $true }
    @{ Cell = 'E19'; Value = '  +3.64%  '; IsPrice = $false }
    @{ Cell = 'D20'; Value = '16.23'; IsPrice = $true }
    @{ Cell = 'E20'; Value = '  -1.46%  '; IsPrice = $false }
    @{ Cell = 'D21'; Value = '488.37'; IsPrice = $true }
    @{ Cell = 'E21'; Value = '  +5.00%  '; IsPrice = $false }
    @{ Cell = 'D22'; Value = '0.720'; IsPrice = $true }
    @{ Cell = 'E22'; Value = '  +1.61%  '; IsPrice = $false }
    @{ Cell = 'D23'; Value = '7.74'; IsPrice = $true }
    @{ Cell = 'E23'; Value = '  +4.85%  '; IsPrice = $false }
    @{ Cell = 'B24'; Value = 'InternetComputer(DFINITY)'; IsPrice = $false }
    @{ Cell = 'C24'; Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'; IsPrice = $false }
    @{ Cell = 'D24'; Value = '13.42'; IsPrice = $true }
    @{ Cell = 'E24'; Value = '  +5.23%  '; IsPrice = $false }
    @{ Cell = 'B25'; Value = 'Litecoin'; IsPrice = $false }
    @{ Cell = 'C25'; Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'; IsPrice = $false }
    @{ Cell = 'D25'; Value = '84.49'; IsPrice = $true }
    @{ Cell = 'E25'; Value = '  +1.80%  '; IsPrice = $false }
    @{ Cell = 'E26'; Value = '  +3.51%  '; IsPrice = $false }
    @{ Cell = 'D27'; Value = '10.11'; IsPrice = $true }
    @{ Cell = 'E27'; Value = '  +0.63%  '; IsPrice = $false }
    @{ Cell = 'E28'; Value = '  +0.05%  '; IsPrice = $false }
    @{ Cell = 'D29'; Value = '8.05'; IsPrice = $true }
    @{ Cell = 'E29'; Value = '  -2.07%  '; IsPrice = $false }
    @{ Cell = 'E30'; Value = '  -1.63%  '; IsPrice = $false }
    @{ Cell = 'D31'; Value = '2.70'; IsPrice = $true }
    @{ Cell = 'E31'; Value = '  +2.55%  '; IsPrice = $false }
    @{ Cell = 'B32'; Value = 'EthereumClassic'; IsPrice = $false }
    @{ Cell = 'C32'; Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'; IsPrice = $false }
    @{ Cell = 'D32'; Value = '29.06'; IsPrice = $true }
    @{ Cell = 'E32'; Value = '  +2.53%  '; IsPrice = $false }
    @{ Cell = 'B33'; Value = 'PEPE'; IsPrice = $false }
    @{ Cell = 'C33'; Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'; IsPrice = $false }
    @{ Cell = 'D33'; Value = '0.0000101'; IsPrice = $true }
    @{ Cell = 'E33'; Value = '  +0.89%  '; IsPrice = $false }
    @{ Cell = 'E34'; Value = '  -3.06%  '; IsPrice = $false }
    @{ Cell = 'E35'; Value = '  +0.06%  '; IsPrice = $false }
    @{ Cell = 'D36'; Value = '5.98'; IsPrice = $true }
    @{ Cell = 'E36'; Value = '  +1.33%  '; IsPrice = $false }
    @{ Cell = 'D37'; Value = '0.993'; IsPrice = $true }
    @{ Cell = 'E37'; Value = '  -0.01%  '; IsPrice = $false }
    @{ Cell = 'D38'; Value = '47.51'; IsPrice = $true }
    @{ Cell = 'E38'; Value = '  -2.21%  '; IsPrice = $false }
    @{ Cell = 'D39'; Value = '2.13'; IsPrice = $true }
    @{ Cell = 'E39'; Value = '  +3.41%  '; IsPrice = $false }
    @{ Cell = 'D40'; Value = '50.23'; IsPrice = $true }
    @{ Cell = 'E40'; Value = '  +1.66%  '; IsPrice = $false }
    @{ Cell = 'D41'; Value = '0.316'; IsPrice = $true }
    @{ Cell = 'E41'; Value = '  +1.21%  '; IsPrice = $false }
    @{ Cell = 'D42'; Value = '0.124'; IsPrice = $true }
    @{ Cell = 'E42'; Value = '  +2.38%  '; IsPrice = $false }
    @{ Cell = 'D43'; Value = '8.69'; IsPrice = $true }
    @{ Cell = 'E43'; Value = '  +0.94%  '; IsPrice = $false }
    @{ Cell = 'D44'; Value = '2.83'; IsPrice = $true }
    @{ Cell = 'E44'; Value = '  -0.54%  '; IsPrice = $false }
    @{ Cell = 'D45'; Value = '2.864.62'; IsPrice = $true }
    @{ Cell = 'E45'; Value = '  +5.68%  '; IsPrice = $false }
    @{ Cell = 'D46'; Value = '389.84'; IsPrice = $true }
    @{ Cell = 'E46'; Value = '  +2.48%  '; IsPrice = $false }
    @{ Cell = 'D47'; Value = '0.0360'; IsPrice = $true }
    @{ Cell = 'E47'; Value = '  +0.01%  '; IsPrice = $false }
    @{ Cell = 'D48'; Value = '136.50'; IsPrice = $true }
    @{ Cell = 'E48'; Value = '  +1.09%  '; IsPrice = $false }
    @{ Cell = 'E49'; Value = '  +0.01%  '; IsPrice = $false }
    @{ Cell = 'D50'; Value = '25.27'; IsPrice = $true }
    @{ Cell = 'E50'; Value = '  +2.83%  '; IsPrice = $false }
    @{ Cell = 'D51'; Value = '2.24'; IsPrice = $true }
    @{ Cell = 'E51'; Value = '  +0.39%  '; IsPrice = $false }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.IsPrice) {
        $rng.Value = "'" + $u.Value
        $rng.Style = "Normal"
    } else {
        $rng.Value = $u.Value
    }
}

Write-Host "Updated $($updates.Count) cells on $($ws.Name)"
